$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 23 new rows (24-46) for the new "flag game" dialogue scene ---
$ws.Range("A24:A46").Insert()

# --- Step 2: fill the new dialogue rows in column A ---
$ws.Cells.Item(24, 1).Value2 = '\n<リリー>バカじゃないの？いい加減にしてよ！'
$ws.Cells.Item(25, 1).Value2 = '\n<シィナ>はぁ～？そっちが悪いにゃ！'
$ws.Cells.Item(26, 1).Value2 = '\n<ライム>はいはいはいはい。落ち着いてー。
どっちが悪いかゲームで決めよー。'
$ws.Cells.Item(27, 1).Value2 = '\n<シィナ>あ？ゲーム？
何にゃ。'
$ws.Cells.Item(28, 1).Value2 = '\n<ライム>旗上げゲーム！'
$ws.Cells.Item(29, 1).Value2 = '\n<リリー>やらないわよバカバカしい。'
$ws.Cells.Item(30, 1).Value2 = '\n<シィナ>負けるのが怖いにゃ？
アタシの瞬発力に勝てるわけないにゃんねー！
しゅっ！しゅっ！しゅっ！'
$ws.Cells.Item(31, 1).Value2 = '\n<リリー>あんた脳みそにアイロンでもかけたの？
シワが無いみたいだけど。'
$ws.Cells.Item(32, 1).Value2 = '\n<シィナ>言ってる意味が分からんにゃーん。
ビビってんのか？
しゅっ！しゅっ！'
$ws.Cells.Item(33, 1).Value2 = '\n<リリー>旗上げは瞬発力より判断力よ。
あんたなんかに負けるわけないでしょ。'
$ws.Cells.Item(34, 1).Value2 = '\n<シィナ>お、やんのか？
どっちが上か、赤白付けてやろうにゃん。'
$ws.Cells.Item(35, 1).Value2 = '\n<ライム>負けた方が謝ること。言い訳無し。
一発勝負。
いくよー？'
$ws.Cells.Item(36, 1).Value2 = '\n<リリー>負けたらクツでもケツでも舐めてやるわよ。'
$ws.Cells.Item(37, 1).Value2 = '\n<シィナ>クツもケツも舐めさせてやるし。'
$ws.Cells.Item(38, 1).Value2 = '\n<ライム>赤\.\.上げ\.ない♪\.\.\.\.白\.\.上げ\.ない♪\.\.\.\.
赤\.\.上げ\.ない\.で♪\.\.\.\.白\.\.上げ\.ない♪'
$ws.Cells.Item(39, 1).Value2 = '\n<ライム>赤\.\.下げ\.て♪\.\.\.\.白\.\.下げ\.て♪\.\.\.\.
白\.\.上げ\.ない\.で♪\.\.\.\.白\.\.下げ\.て♪'
$ws.Cells.Item(40, 1).Value2 = '\n<ライム>赤\.\.上げ\.ない♪\.\.\.\.白\.\.上げ\.ない♪\.\.\.\.
白\.\.上げ\.ない\.で\.\^'
$ws.Cells.Item(41, 1).Value2 = '\n<シィナ>上げさせろにゃ！'
$ws.Cells.Item(42, 1).Value2 = '\n<リリー>あーもうムカツクわね！'
$ws.Cells.Item(43, 1).Value2 = '\n<ライム>ちょ、ちょっとまってよー！
これからなのにー！'
$ws.Cells.Item(44, 1).Value2 = '\n<研究員1>こいつらほんと仲良いっすね。'
$ws.Cells.Item(45, 1).Value2 = '\n<研究員2>どっちかが勝ったらまた喧嘩になっちゃいますからね。
ライムちゃん、いい子ですね。'
$ws.Cells.Item(46, 1).Value2 = '\n<研究員1>いや、多分あれ素っすよ。'

# --- Step 3: mirror column A into column B for the already-existing dialogue rows (1-23) ---
for ($i = 1; $i -le 23; $i++) {
    $ws.Cells.Item($i, 2).Value2 = $ws.Cells.Item($i, 1).Value2
}

# --- Step 4: mirror column A into column B for the event-name rows (now 47-49) ---
for ($i = 47; $i -le 49; $i++) {
    $ws.Cells.Item($i, 2).Value2 = $ws.Cells.Item($i, 1).Value2
}

# --- Step 5: row 50 ("オーク2 ") becomes the new "group leader" entry; the original label
# --- moves into column B, column A gets the new "オーク2 リーダ" label ---
$ws.Cells.Item(50, 2).Value2 = $ws.Cells.Item(50, 1).Value2
$ws.Cells.Item(50, 1).Value2 = 'オーク2 リーダ'

# --- Step 6: mirror column A into column B for the remaining event-name rows (51-52) ---
for ($i = 51; $i -le 52; $i++) {
    $ws.Cells.Item($i, 2).Value2 = $ws.Cells.Item($i, 1).Value2
}

# --- Step 7: append two new rows at the end for the flag-game props ---
$ws.Cells.Item(53, 1).Value2 = '旗'
$ws.Cells.Item(54, 1).Value2 = '旗２'
